# Apply the cryptos-list refresh for Sun Oct 15 22:11:14 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:E hold text (coin name / link / price / % change) even when the
# string looks numeric (e.g. "210.44", "1.00"). Force the cells we touch to
# stay text-formatted first, so Excel does not silently coerce them to numbers
# and drop formatting such as trailing zeros.
$changedCells = @("D2","E2","D3","E3","E4","D5","E5","E6","E7","D8","E8","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","D27","E27","E28","E29","E30","E31","D32","E32","D33","E33","D34","E34","E35","D36","E36","D37","E37","E38","D39","E39","E40","B41","C41","D41","E41","B42","C42","D42","B43","C43","D43","E43","D44","E44","D45","E45","E46","D47","E47","D48","E48","D49","E49","B50","C50","D50","E50","B51","C51","D51","E51")
foreach ($addr in $changedCells) { $ws.Range($addr).NumberFormat = "@" }

# Row 2
$ws.Range("D2").Value = "27.264.37"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").Value = "1.565.56"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.40%  "

# Row 5
$ws.Range("D5").Value = "210.44"
$ws.Range("E5").Value = "  +1.35%  "

# Row 6
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("E7").Value = "  -0.58%  "

# Row 8
$ws.Range("D8").Value = "22.14"
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("D10").Value = "0.0597"
$ws.Range("E10").Value = "  -0.67%  "

# Row 11
$ws.Range("D11").Value = "0.0869"
$ws.Range("E11").Value = "  +1.50%  "

# Row 12
$ws.Range("D12").Value = "1.784.88"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").Value = "1.572.05"
$ws.Range("E13").Value = "  +0.47%  "

# Row 14
$ws.Range("D14").Value = "3.78"
$ws.Range("E14").Value = "  +0.47%  "

# Row 15
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -0.34%  "

# Row 16
$ws.Range("D16").Value = "27.225.52"
$ws.Range("E16").Value = "  +0.88%  "

# Row 17
$ws.Range("D17").Value = "61.92"
$ws.Range("E17").Value = "  -0.24%  "

# Row 18
$ws.Range("D18").Value = "7.48"
$ws.Range("E18").Value = "  +1.84%  "

# Row 19
$ws.Range("D19").Value = "217.27"
$ws.Range("E19").Value = "  +0.22%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0702"
$ws.Range("E20").Value = "  -0.49%  "

# Row 21
$ws.Range("E21").Value = "  -0.43%  "

# Row 22
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  +0.53%  "

# Row 23
$ws.Range("D23").Value = "9.22"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("D25").Value = "152.91"
$ws.Range("E25").Value = "  -0.12%  "

# Row 26
$ws.Range("D26").Value = "6.63"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("D27").Value = "15.02"
$ws.Range("E27").Value = "  -0.52%  "

# Row 28
$ws.Range("E28").Value = "  +1.46%  "

# Row 29
$ws.Range("E29").Value = "  -0.44%  "

# Row 30
$ws.Range("E30").Value = "  +1.78%  "

# Row 31
$ws.Range("E31").Value = "  -0.23%  "

# Row 32
$ws.Range("D32").Value = "3.24"
$ws.Range("E32").Value = "  -0.04%  "

# Row 33
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  +1.27%  "

# Row 34
$ws.Range("D34").Value = "1.435.34"
$ws.Range("E34").Value = "  +0.75%  "

# Row 35
$ws.Range("E35").Value = "  +3.66%  "

# Row 36
$ws.Range("D36").Value = "1.61"
$ws.Range("E36").Value = "  -0.35%  "

# Row 37
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -0.23%  "

# Row 38
$ws.Range("E38").Value = "  +0.52%  "

# Row 39
$ws.Range("D39").Value = "0.531"
$ws.Range("E39").Value = "  -0.75%  "

# Row 40
$ws.Range("E40").Value = "  +2.09%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  -0.30%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.00"

# Row 43
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  +0.26%  "

# Row 44
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -1.07%  "

# Row 45
$ws.Range("D45").Value = "64.50"
$ws.Range("E45").Value = "  -0.50%  "

# Row 46
$ws.Range("E46").Value = "  -1.01%  "

# Row 47
$ws.Range("D47").Value = "1.699.57"
$ws.Range("E47").Value = "  -0.14%  "

# Row 48
$ws.Range("D48").Value = "85.99"
$ws.Range("E48").Value = "  -1.71%  "

# Row 49
$ws.Range("D49").Value = "0.0525"
$ws.Range("E49").Value = "  +1.10%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  -0.85%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0952"
$ws.Range("E51").Value = "  -0.96%  "
